$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Fill in the new Sprint #2 row of data (row 5)
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = "A, B, D, E, F, I, K, N"
$ws.Range("E5").Value = "Multiple Team meetings, along with Conference call with TARDEC, established basic website template, User stories for frontend"
$ws.Range("F5").Value = "userstories_master, SkyPIsite, ui_user_stories"
$ws.Range("G5").Value = "ui_user_stories.txt, userstories_master.txt, SkyPIsite directory"
$ws.Range("H5").Value = "Both User Stories files frontend sections, SkyPIsite folder was our first template but has recently been updated with a newer template (view commits in GitHub for more details)"

# Wrap text on the newly filled description cells and resize the row to fit
$ws.Range("E5:H5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 94.5

# Move the active selection like the author left it
[void]$ws.Range("G6").Select()
